# Update column G ("K") values in Sheet1 with new data, per the
# "regen save_data to use K instead of Strike#, regen std/mean,
# calc and write s_vals" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newK = @{
    2  = 6
    3  = 4
    4  = 0
    5  = 3
    6  = 1
    7  = 2
    8  = 4
    9  = 8
    10 = 4
    11 = 3
    12 = 1
    13 = 2
    14 = 1
    15 = 5
    16 = 7
    17 = 3
    18 = 1
    19 = 4
    20 = 4
    21 = 3
    22 = 7
    23 = 2
    24 = 1
    25 = 4
    26 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
